$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.338.11'
$ws.Range("E2").Value = '  +6.04%  '
$ws.Range("D3").Value = '2.365.44'
$ws.Range("E3").Value = '  +2.20%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '109.54'
$ws.Range("E5").Value = '  +2.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '310.12'
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  +1.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.28'
$ws.Range("E10").Value = '  +2.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0918'
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.50'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.988'
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").Value = '2.724.89'
$ws.Range("E15").Value = '  +2.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.41'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").Value = '2.373.39'
$ws.Range("E17").Value = '  +2.39%  '
$ws.Range("D18").Value = '45.331.93'
$ws.Range("E18").Value = '  +6.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.34'
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.93'
$ws.Range("E20").Value = '  +6.25%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000107'
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.41'
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.49'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '259.78'
$ws.Range("E24").Value = '  -2.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("E25").Value = '  +3.24%  '
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.18'
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("E28").Value = '  -5.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.37'
$ws.Range("E29").Value = '  +3.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0972'
$ws.Range("E30").Value = '  +10.91%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.37'
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.99'
$ws.Range("E32").Value = '  -1.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '169.37'
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("E34").Value = '  +6.35%  '
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("E36").Value = '  +4.60%  '
$ws.Range("E37").Value = '  +2.98%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.96'
$ws.Range("E38").Value = '  +7.67%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'
$ws.Range("E39").Value = '  +4.48%  '
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.75'
$ws.Range("E41").Value = '  +8.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.32'
$ws.Range("E42").Value = '  -4.16%  '
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '69.88'
$ws.Range("E44").Value = '  -1.48%  '
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '83.02'
$ws.Range("E47").Value = '  +7.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '112.58'
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.50'
$ws.Range("E49").Value = '  +5.02%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.676.29'
$ws.Range("E50").Value = '  +1.30%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.13'
$ws.Range("E51").Value = '  +3.79%  '
